# The commit inserts one new weekly price record for "Papa" (potato) at
# what becomes row 55 of the sheet. Every row below the insertion point
# shifts down by one (old row 55 -> new row 56, ..., old row 106 -> new
# row 107), and the sheet's used range grows from A1:R106 to A1:R107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 55..106 down one row, opening up a blank row 55.
$ws.Rows.Item(55).Insert()

# Populate the newly opened row 55 with the new record.
$ws.Cells.Item(55, 1).Value  = 11
$ws.Cells.Item(55, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value  = "Bíobío"
$ws.Cells.Item(55, 4).Value  = 44447
$ws.Cells.Item(55, 5).Value  = 8
$ws.Cells.Item(55, 6).Value  = 100114001
$ws.Cells.Item(55, 7).Value  = "Papa"
$ws.Cells.Item(55, 8).Value  = "Patagonia"
$ws.Cells.Item(55, 9).Value  = "1a (guarda)"
$ws.Cells.Item(55, 10).Value = 2000
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 8500
$ws.Cells.Item(55, 13).Value = 8250
$ws.Cells.Item(55, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(55, 16).Value = 330
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"
